$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price values: force text so Excel does not reinterpret
# numeric-looking strings (e.g. "628.10", "0.0000247") as numbers,
# which would drop trailing/insignificant zeros and change the display.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "91.316.38"
$ws.Range("E2").Value = "  +1.80%  "

Set-TextValue $ws.Range("D3") "3.177.90"
$ws.Range("E3").Value = "  +4.97%  "

Set-TextValue $ws.Range("D5") "216.09"
$ws.Range("E5").Value = "  +2.59%  "

Set-TextValue $ws.Range("D6") "628.10"
$ws.Range("E6").Value = "  +2.69%  "

$ws.Range("E7").Value = "  +31.99%  "

$ws.Range("E8").Value = "  +2.64%  "

Set-TextValue $ws.Range("D10") "3.174.59"
$ws.Range("E10").Value = "  +4.97%  "

Set-TextValue $ws.Range("D11") "0.768"
$ws.Range("E11").Value = "  +15.54%  "

Set-TextValue $ws.Range("D12") "0.203"
$ws.Range("E12").Value = "  +7.65%  "

Set-TextValue $ws.Range("D13") "0.0000247"
$ws.Range("E13").Value = "  +2.90%  "

Set-TextValue $ws.Range("D14") "5.68"
$ws.Range("E14").Value = "  +6.55%  "

Set-TextValue $ws.Range("D15") "35.14"
$ws.Range("E15").Value = "  +9.40%  "

Set-TextValue $ws.Range("D16") "90.969.25"
$ws.Range("E16").Value = "  +1.84%  "

Set-TextValue $ws.Range("D17") "3.759.99"
$ws.Range("E17").Value = "  +4.77%  "

Set-TextValue $ws.Range("D18") "3.154.88"
$ws.Range("E18").Value = "  +3.79%  "

Set-TextValue $ws.Range("D19") "3.72"
$ws.Range("E19").Value = "  +13.12%  "

Set-TextValue $ws.Range("D20") "14.61"
$ws.Range("E20").Value = "  +9.31%  "

Set-TextValue $ws.Range("D21") "471.89"
$ws.Range("E21").Value = "  +11.69%  "

$ws.Range("E22").Value = "  -2.53%  "

Set-TextValue $ws.Range("D23") "9.20"
$ws.Range("E23").Value = "  +11.89%  "

Set-TextValue $ws.Range("D24") "5.18"
$ws.Range("E24").Value = "  +3.11%  "

Set-TextValue $ws.Range("D25") "96.58"
$ws.Range("E25").Value = "  +17.38%  "

Set-TextValue $ws.Range("D26") "5.95"
$ws.Range("E26").Value = "  +11.27%  "

Set-TextValue $ws.Range("D27") "12.49"
$ws.Range("E27").Value = "  +8.43%  "

Set-TextValue $ws.Range("D28") "3.337.08"
$ws.Range("E28").Value = "  +4.52%  "

Set-TextValue $ws.Range("D30") "9.40"
$ws.Range("E30").Value = "  +13.60%  "

$ws.Range("E31").Value = "  +0.78%  "

$ws.Range("E32").Value = "  +0.08%  "

Set-TextValue $ws.Range("D33") "27.61"
$ws.Range("E33").Value = "  +21.83%  "

Set-TextValue $ws.Range("D34") "528.61"
$ws.Range("E34").Value = "  +5.75%  "

Set-TextValue $ws.Range("D35") "0.193"
$ws.Range("E35").Value = "  +43.67%  "

$ws.Range("E36").Value = "  +8.07%  "

$ws.Range("E37").Value = "  -1.89%  "

$ws.Range("E38").Value = "  +10.13%  "

$ws.Range("E39").Value = "  +5.68%  "

$ws.Range("E40").Value = "  +6.19%  "

Set-TextValue $ws.Range("D41") "0.0872"
$ws.Range("E41").Value = "  +25.90%  "

Set-TextValue $ws.Range("D42") "22.28"
$ws.Range("E42").Value = "  +0.14%  "

Set-TextValue $ws.Range("D43") "0.423"
$ws.Range("E43").Value = "  +17.88%  "

$ws.Range("E44").Value = "  +0.02%  "

$ws.Range("E45").Value = "  +9.54%  "

Set-TextValue $ws.Range("D47") "0.709"
$ws.Range("E47").Value = "  +20.72%  "

Set-TextValue $ws.Range("D48") "152.38"
$ws.Range("E48").Value = "  +5.34%  "

Set-TextValue $ws.Range("D49") "4.66"
$ws.Range("E49").Value = "  +11.67%  "

$ws.Range("E50").Value = "  +12.54%  "

Set-TextValue $ws.Range("D51") "45.40"
$ws.Range("E51").Value = "  +4.42%  "
